$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (changed) date column C for all data rows (2-505)
#    from 45192 to 45202.
$ws.Range("C2:C505").Value = 45202

# 2. Rows 6 and 7 swap identities: what used to be row 6 (A 59850-2022) moves
#    to row 7, and what used to be row 7 (A 11067-2023) moves to row 6 - with
#    updated figures for "A 11067-2023" (Signalarter I6 becomes 3, Alla arter
#    Q6 becomes 6, and a new species "Bollvitmossa" is added to the species
#    list).

# New row 6: A 11067-2023
$ws.Range("A6").Value = "A 11067-2023"
$ws.Range("B6").Value = 44992
$ws.Range("C6").Value = 45202
$ws.Range("D6").Value = "DALARNAS LÄN"
$ws.Range("E6").Value = "LEKSAND"
$ws.Range("G6").Value = 2.9
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 3
$ws.Range("J6").Value = 3
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 3
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 6
$ws.Range("R6").Value = "Reliktbock`r`nRosenticka`r`nUllticka`r`nBollvitmossa`r`nStuplav`r`nVedticka"
$ws.Range("S6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LEKSAND/artfynd/A 11067-2023.xlsx", "A 11067-2023")'
$ws.Range("T6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LEKSAND/kartor/A 11067-2023.png", "A 11067-2023")'
$ws.Range("V6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LEKSAND/klagomål/A 11067-2023.docx", "A 11067-2023")'
$ws.Range("W6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LEKSAND/klagomålsmail/A 11067-2023.docx", "A 11067-2023")'
$ws.Range("X6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LEKSAND/tillsyn/A 11067-2023.docx", "A 11067-2023")'
$ws.Range("Y6").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LEKSAND/tillsynsmail/A 11067-2023.docx", "A 11067-2023")'

# New row 7: A 59850-2022
$ws.Range("A7").Value = "A 59850-2022"
$ws.Range("B7").Value = 44901
$ws.Range("C7").Value = 45202
$ws.Range("D7").Value = "DALARNAS LÄN"
$ws.Range("E7").Value = "LEKSAND"
$ws.Range("G7").Value = 1.7
$ws.Range("H7").Value = 3
$ws.Range("I7").Value = 1
$ws.Range("J7").Value = 2
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 2
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 5
$ws.Range("R7").Value = "Tretåig hackspett`r`nViolettgrå tagellav`r`nSvavelriska`r`nBlåsippa`r`nRevlummer"
$ws.Range("S7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LEKSAND/artfynd/A 59850-2022.xlsx", "A 59850-2022")'
$ws.Range("T7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LEKSAND/kartor/A 59850-2022.png", "A 59850-2022")'
$ws.Range("V7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LEKSAND/klagomål/A 59850-2022.docx", "A 59850-2022")'
$ws.Range("W7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LEKSAND/klagomålsmail/A 59850-2022.docx", "A 59850-2022")'
$ws.Range("X7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LEKSAND/tillsyn/A 59850-2022.docx", "A 59850-2022")'
$ws.Range("Y7").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LEKSAND/tillsynsmail/A 59850-2022.docx", "A 59850-2022")'

# Keep the row heights the same as the rest of the sheet (15), since the
# wrapped-text autosize would otherwise change it after editing the R
# column content.
$ws.Rows.Item(6).RowHeight = 15
$ws.Rows.Item(7).RowHeight = 15
